# Generate Report for Handoff
# Update status + timestamps to reflect that translations are now
# "Ready for handoff" instead of "In Translation", and bump the
# latest handoff / xliff-generation timestamps accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-04 13:03:10"

# zh-cn sheet: "Status" column (C2) and "Latest Handoff Datetime" (H2).
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-04 13:03:01"

# de-de sheet: "Status" column (C2) and "Latest Handoff Datetime" (H2).
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-04 13:03:10"

# The status-text columns grew wider ("In Translation" -> "Ready for
# handoff"), so widen the columns that hold that text on every sheet to
# match Excel's automatic column-width recalculation (target stored
# width ~= 17.216 characters).
$newStatusColWidth = 16.3333333333333
$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$dede.Columns.Item(3).ColumnWidth = $newStatusColWidth
